# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates described in the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 189.4
$ws.Range("I8").Value = 56.285713
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 168.857139
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = -29.85713900000002
$ws.Range("N8").Value = -1778
$ws.Range("H18").Value = 1273.7693
$ws.Range("J18").Value = 525.5
$ws.Range("L18").Value = 525.5
$ws.Range("N18").Value = -1093.5
$ws.Range("H116").Value = 4473.4
$ws.Range("I116").Value = 4450
$ws.Range("J116").Value = 4489
$ws.Range("K116").Value = 4450
$ws.Range("L116").Value = 4489
$ws.Range("M116").Value = -1008
$ws.Range("N116").Value = -11373
$ws.Range("H135").Value = 786.53845
$ws.Range("I135").Value = 692.5
$ws.Range("K135").Value = 6232.5
$ws.Range("M135").Value = -3697.5
$ws.Range("H137").Value = 4634255.5
$ws.Range("I137").Value = 7145496.5
$ws.Range("K137").Value = 21436489.5
$ws.Range("M137").Value = -21433939.5
$ws.Range("H138").Value = 2963.4656
$ws.Range("J138").Value = 4014.5334
$ws.Range("L138").Value = 12043.6002
$ws.Range("N138").Value = -22323.6002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 16022.833
$ws.Range("I26").Value = 2380.3333
$ws.Range("J26").Value = 29665.334
$ws.Range("K26").Value = 2380.3333
$ws.Range("L26").Value = 29665.334
$ws.Range("M26").Value = -2050.3333
$ws.Range("N26").Value = -30325.334
$ws.Range("H32").Value = 28366.227
$ws.Range("I32").Value = 29407.477
$ws.Range("K32").Value = 29407.477
$ws.Range("M32").Value = -29120.477
$ws.Range("H40").Value = 42500
$ws.Range("J40").Value = 60000
$ws.Range("L40").Value = 60000
$ws.Range("N40").Value = -60352
$ws.Range("H64").Value = 5017520.5
$ws.Range("I64").Value = 3356694
$ws.Range("K64").Value = 3356694
$ws.Range("M64").Value = -3356446
$ws.Range("H67").Value = 5017520.5
$ws.Range("I67").Value = 3356694
$ws.Range("K67").Value = 3356694
$ws.Range("M67").Value = -3355836
$ws.Range("H74").Value = 1673639.6
$ws.Range("I74").Value = 5002012
$ws.Range("J74").Value = 9453.5
$ws.Range("K74").Value = 5002012
$ws.Range("L74").Value = 9453.5
$ws.Range("M74").Value = -5001138
$ws.Range("N74").Value = -11201.5
$ws.Range("H77").Value = 1673639.6
$ws.Range("I77").Value = 5002012
$ws.Range("J77").Value = 9453.5
$ws.Range("K77").Value = 25010060
$ws.Range("L77").Value = 47267.5
$ws.Range("M77").Value = -25005692
$ws.Range("N77").Value = -56003.5
$ws.Range("H132").Value = 5274.9355
$ws.Range("I132").Value = 3715.9167
$ws.Range("K132").Value = 11147.7501
$ws.Range("M132").Value = -8617.750100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 126174.19
$ws.Range("I86").Value = 1018.7
$ws.Range("K86").Value = 1018.7
$ws.Range("M86").Value = 104.3
$ws.Range("H89").Value = 126174.19
$ws.Range("I89").Value = 1018.7
$ws.Range("K89").Value = 5093.5
$ws.Range("M89").Value = 522.5
$ws.Range("H92").Value = 62800.332
$ws.Range("J92").Value = 62800.332
$ws.Range("L92").Value = 62800.332
$ws.Range("N92").Value = -67792.33199999999
$ws.Range("H94").Value = 1086.75
$ws.Range("I94").Value = 1062.8182
$ws.Range("J94").Value = 1350
$ws.Range("K94").Value = 1062.8182
$ws.Range("L94").Value = 1350
$ws.Range("M94").Value = -611.8181999999999
$ws.Range("N94").Value = -2252
$ws.Range("H107").Value = 1211.3846
$ws.Range("I107").Value = 1068.45
$ws.Range("K107").Value = 1068.45
$ws.Range("M107").Value = 851.55
$ws.Range("H134").Value = 3715.3447
$ws.Range("I134").Value = 1924
$ws.Range("J134").Value = 10582.167
$ws.Range("K134").Value = 5772
$ws.Range("L134").Value = 31746.501
$ws.Range("M134").Value = -3237
$ws.Range("N134").Value = -36816.501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30307918
$ws.Range("I31").Value = 55558964
$ws.Range("J31").Value = 6665.6
$ws.Range("K31").Value = 55558964
$ws.Range("L31").Value = 6665.6
$ws.Range("M31").Value = -55558669
$ws.Range("N31").Value = -7255.6
$ws.Range("H34").Value = 30307918
$ws.Range("I34").Value = 55558964
$ws.Range("J34").Value = 6665.6
$ws.Range("K34").Value = 55558964
$ws.Range("L34").Value = 6665.6
$ws.Range("M34").Value = -55558762
$ws.Range("N34").Value = -7069.6
$ws.Range("H58").Value = 4114.885
$ws.Range("I58").Value = 3038.6667
$ws.Range("K58").Value = 3038.6667
$ws.Range("M58").Value = -2835.6667
$ws.Range("H103").Value = 36248.25
$ws.Range("J103").Value = 79993
$ws.Range("L103").Value = 79993
$ws.Range("N103").Value = -82337
$ws.Range("H132").Value = 3556
$ws.Range("I132").Value = 3020.9412
$ws.Range("K132").Value = 9062.8236
$ws.Range("M132").Value = -6532.8236
$ws.Range("H136").Value = 4114.885
$ws.Range("I136").Value = 3038.6667
$ws.Range("K136").Value = 9116.000100000001
$ws.Range("M136").Value = -6566.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 800.0714
$ws.Range("I5").Value = 428.85715
$ws.Range("J5").Value = 1171.2858
$ws.Range("K5").Value = 1286.57145
$ws.Range("L5").Value = 3513.8574
$ws.Range("M5").Value = -1174.57145
$ws.Range("N5").Value = -3737.8574
$ws.Range("H113").Value = 2132.238
$ws.Range("I113").Value = 3045.25
$ws.Range("K113").Value = 9135.75
$ws.Range("M113").Value = -6965.75
$ws.Range("H122").Value = 23673
$ws.Range("I122").Value = 35109.5
$ws.Range("K122").Value = 315985.5
$ws.Range("M122").Value = -313535.5
$ws.Range("H129").Value = 18521468
$ws.Range("J129").Value = 41668010
$ws.Range("L129").Value = 125004030
$ws.Range("N129").Value = -125014030
$ws.Range("H135").Value = 800.0714
$ws.Range("I135").Value = 428.85715
$ws.Range("J135").Value = 1171.2858
$ws.Range("K135").Value = 3859.71435
$ws.Range("L135").Value = 10541.5722
$ws.Range("M135").Value = -1324.71435
$ws.Range("N135").Value = -15611.5722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 9060.5
$ws.Range("I19").Value = 336.66666
$ws.Range("J19").Value = 12799.286
$ws.Range("K19").Value = 336.66666
$ws.Range("L19").Value = 12799.286
$ws.Range("M19").Value = -48.66665999999998
$ws.Range("N19").Value = -13375.286
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H54").Value = 5555
$ws.Range("J54").Value = 5555
$ws.Range("L54").Value = 5555
$ws.Range("N54").Value = -6335
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H132").Value = 3452.6316
$ws.Range("I132").Value = 1576.8
$ws.Range("K132").Value = 4730.4
$ws.Range("M132").Value = -2200.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -888
$ws.Range("H34").Value = 5000
$ws.Range("I34").Value = 5000
$ws.Range("K34").Value = 5000
$ws.Range("M34").Value = -4828
$ws.Range("H46").Value = 4650.727
$ws.Range("I46").Value = 909.6
$ws.Range("K46").Value = 909.6
$ws.Range("M46").Value = -721.6
$ws.Range("H61").Value = 2255.1538
$ws.Range("I61").Value = 2228.9092
$ws.Range("J61").Value = 2399.5
$ws.Range("K61").Value = 2228.9092
$ws.Range("L61").Value = 2399.5
$ws.Range("M61").Value = -2026.9092
$ws.Range("N61").Value = -2803.5
$ws.Range("H68").Value = 2828.5881
$ws.Range("I68").Value = 2139.0667
$ws.Range("J68").Value = 8000
$ws.Range("K68").Value = 2139.0667
$ws.Range("L68").Value = 8000
$ws.Range("M68").Value = -1390.0667
$ws.Range("N68").Value = -9498
$ws.Range("H71").Value = 2828.5881
$ws.Range("I71").Value = 2139.0667
$ws.Range("J71").Value = 8000
$ws.Range("K71").Value = 10695.3335
$ws.Range("L71").Value = 40000
$ws.Range("M71").Value = -6951.333499999999
$ws.Range("N71").Value = -47488
$ws.Range("H113").Value = 2255.1538
$ws.Range("I113").Value = 2228.9092
$ws.Range("J113").Value = 2399.5
$ws.Range("K113").Value = 2228.9092
$ws.Range("L113").Value = 2399.5
$ws.Range("M113").Value = -58.90920000000006
$ws.Range("N113").Value = -6739.5
$ws.Range("H132").Value = 6331.9
$ws.Range("I132").Value = 3815.9333
$ws.Range("K132").Value = 11447.7999
$ws.Range("M132").Value = -8917.7999
$ws.Range("H136").Value = 6563.4375
$ws.Range("I136").Value = 5301.6665
$ws.Range("K136").Value = 15904.9995
$ws.Range("M136").Value = -13354.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 66674
$ws.Range("J43").Value = 79995
$ws.Range("L43").Value = 79995
$ws.Range("N43").Value = -80293
$ws.Range("H113").Value = 411.86667
$ws.Range("I113").Value = 411.86667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1235.60001
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 934.3999899999999
$ws.Range("N113").ClearContents()
$ws.Range("H124").Value = 114996.336
$ws.Range("J124").Value = 114996.336
$ws.Range("L124").Value = 114996.336
$ws.Range("N124").Value = -124816.336
$ws.Range("H136").Value = 1876.475
$ws.Range("I136").Value = 1081.1143
$ws.Range("J136").Value = 7444
$ws.Range("K136").Value = 3243.3429
$ws.Range("L136").Value = 22332
$ws.Range("M136").Value = -693.3428999999996
$ws.Range("N136").Value = -27432
